$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range (rows 2-88) before writing the new, reordered dataset
$ws.Range("A2:B88").Clear()

$data = New-Object "object[,]" 75,2
$data[0,0] = "Kreatinine"
$data[0,1] = 10
$data[1,0] = "Length"
$data[1,1] = 10
$data[2,0] = "Vrij T4"
$data[2,1] = 10
$data[3,0] = "25-OH Vitamine D"
$data[3,1] = 10
$data[4,0] = "Neurological"
$data[4,1] = 10
$data[5,0] = "Trombocyten"
$data[5,1] = 10
$data[6,0] = "Kalium"
$data[6,1] = 10
$data[7,0] = "Specialisms_hospitalization"
$data[7,1] = 10
$data[8,0] = "LDH"
$data[8,1] = 10
$data[9,0] = "CKD-EPI eGFR"
$data[9,1] = 10
$data[10,0] = "KIC"
$data[10,1] = 10
$data[11,0] = "Glucose/PL"
$data[11,1] = 10
$data[12,0] = "Psychofarmaca"
$data[12,1] = 10
$data[13,0] = "Alk.Fosf."
$data[13,1] = 10
$data[14,0] = "Icterische index"
$data[14,1] = 10
$data[15,0] = "Lipemische index"
$data[15,1] = 10
$data[16,0] = "RRdiast"
$data[16,1] = 10
$data[17,0] = "Neoplasms"
$data[17,1] = 10
$data[18,0] = "LON"
$data[18,1] = 10
$data[19,0] = "Natrium"
$data[19,1] = 10
$data[20,0] = "ALAT"
$data[20,1] = 10
$data[21,0] = "BMI"
$data[21,1] = 10
$data[22,0] = "LDL-Cholesterol"
$data[22,1] = 10
$data[23,0] = "Radiologic_investigations"
$data[23,1] = 10
$data[24,0] = "Calcium"
$data[24,1] = 10
$data[25,0] = "Musculoskeletal"
$data[25,1] = 9
$data[26,0] = "GYN"
$data[26,1] = 9
$data[27,0] = "FSH"
$data[27,1] = 9
$data[28,0] = "ASAT"
$data[28,1] = 9
$data[29,0] = "Leukocyten"
$data[29,1] = 9
$data[30,0] = "URO"
$data[30,1] = 9
$data[31,0] = "Endocrine/metabolic"
$data[31,1] = 9
$data[32,0] = "NEU"
$data[32,1] = 9
$data[33,0] = "MCV"
$data[33,1] = 9
$data[34,0] = "RDW"
$data[34,1] = 9
$data[35,0] = "END"
$data[35,1] = 9
$data[36,0] = "KIN"
$data[36,1] = 9
$data[37,0] = "Genitourinary"
$data[37,1] = 9
$data[38,0] = "MET"
$data[38,1] = 9
$data[39,0] = "MDL"
$data[39,1] = 9
$data[40,0] = "LOG"
$data[40,1] = 9
$data[41,0] = "Hemolytische index"
$data[41,1] = 8
$data[42,0] = "Tot. Bilirubine"
$data[42,1] = 8
$data[43,0] = "SUM"
$data[43,1] = 8
$data[44,0] = "Pregnancy complications"
$data[44,1] = 8
$data[45,0] = "Total_amount_ICD10s"
$data[45,1] = 8
$data[46,0] = "KLG"
$data[46,1] = 8
$data[47,0] = "KCH"
$data[47,1] = 8
$data[48,0] = "Gamma-GT"
$data[48,1] = 8
$data[49,0] = "RRsyst"
$data[49,1] = 7
$data[50,0] = "Hemoglobine"
$data[50,1] = 7
$data[51,0] = "KEN"
$data[51,1] = 7
$data[52,0] = "Ureum"
$data[52,1] = 7
$data[53,0] = "ANE"
$data[53,1] = 7
$data[54,0] = "KNO"
$data[54,1] = 7
$data[55,0] = "KLZ"
$data[55,1] = 7
$data[56,0] = "NEF"
$data[56,1] = 7
$data[57,0] = "ORT"
$data[57,1] = 6
$data[58,0] = "Albumine"
$data[58,1] = 6
$data[59,0] = "RAD"
$data[59,1] = 6
$data[60,0] = "CAR"
$data[60,1] = 5
$data[61,0] = "AUD"
$data[61,1] = 5
$data[62,0] = "Dermatologic"
$data[62,1] = 5
$data[63,0] = "TSH"
$data[63,1] = 5
$data[64,0] = "HR"
$data[64,1] = 5
$data[65,0] = "Congenital anomalies"
$data[65,1] = 4
$data[66,0] = "RTH"
$data[66,1] = 4
$data[67,0] = "GGZ"
$data[67,1] = 4
$data[68,0] = "Sense organs"
$data[68,1] = 3
$data[69,0] = "HEM"
$data[69,1] = 2
$data[70,0] = "FYS"
$data[70,1] = 2
$data[71,0] = "KNE"
$data[71,1] = 1
$data[72,0] = "NUC"
$data[72,1] = 1
$data[73,0] = "Infectious diseases"
$data[73,1] = 1
$data[74,0] = "Iron-tablets"
$data[74,1] = 1

$ws.Range("A2:B76").Value = $data

# Re-apply the Feature-column style (bold, centered, bordered) that Clear() removed
$ws.Range("A1").Copy()
$ws.Range("A2:A76").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()